$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the "n/a" -> "n.a." typo for the fps rows (D3:D8)
$ws.Range("D3:D8").Value = "n.a."

# Correct the mutation value in G7, which should alternate between
# "Inversion Mutation " and "Arithmetic Mutation " like the other rows
$ws.Range("G7").Value = "Inversion Mutation "

# Update the active selection to match the saved view
$ws.Range("E21").Select()
